# Insert a new weekly record row at row 58 ("Poroto granado" – Macroferia
# Regional de Talca). Excel shifts rows 58..133 down to 59..134 and the
# sheet's used range grows from R133 to R134, matching the rest of the
# diff (every subsequent row's values are simply the previous row's
# values, now one row further down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 58, pushing everything below
# (including the former last row, 133) down by one.
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with the new weekly record.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44638
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = 100112030
$ws.Range("G58").Value = "Poroto granado"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 300
$ws.Range("K58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = 20000
$ws.Range("N58").Value = "$/saco 25 kilos"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 800
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
